$wb = $excel.ActiveWorkbook

# Sheet1: keep data/styling as-is, just move the active selection to B1
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Select() | Out-Null

# Add a second sheet right after Sheet1 and name it "without EndRow"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "without EndRow"

# Row 1: numbers 1..5 in A1:E1, shared string "g" in I1
$ws2.Cells.Item(1,1).Value = 1
$ws2.Cells.Item(1,2).Value = 2
$ws2.Cells.Item(1,3).Value = 3
$ws2.Cells.Item(1,4).Value = 4
$ws2.Cells.Item(1,5).Value = 5
$ws2.Cells.Item(1,9).Value = "g"

# Row 2: same as row 1
$ws2.Cells.Item(2,1).Value = 1
$ws2.Cells.Item(2,2).Value = 2
$ws2.Cells.Item(2,3).Value = 3
$ws2.Cells.Item(2,4).Value = 4
$ws2.Cells.Item(2,5).Value = 5
$ws2.Cells.Item(2,9).Value = "g"

# Row 3: the "#! FINISH" marker + explanation text
$ws2.Cells.Item(3,1).Value = "#! FINISH"
$ws2.Cells.Item(3,2).Value = "<-- this ``#! FINISH`` should be in output, because the line above hasn't ``END_ROW`` that couse finish rendering on out of column limit reached (16384 columns)"

$ws2.Range("B3").Select() | Out-Null
